$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOW2")

# Clear the placeholder "one seed only" labels in B1 and J1 (keep cell style)
$ws.Range("B1").ClearContents()
$ws.Range("J1").ClearContents()

# Left table (TRAINING block): columns B-G, row by row
$ws.Range("B3").Value = "0.6813±0.0024"
$ws.Range("C3").Value = "0.6375±0.0062"
$ws.Range("D3").Value = "0.4664±0.0041"
$ws.Range("E3").Value = "0.3022±0.0071"
$ws.Range("F3").Value = "0.7024±0.0048"
$ws.Range("G3").Value = "0.5387±0.0046"
$ws.Range("B5").Value = "0.655±0.0023"
$ws.Range("C5").Value = "0.5553±0.0038"
$ws.Range("D5").Value = "0.5077±0.0053"
$ws.Range("E5").Value = "0.297±0.0067"
$ws.Range("F5").Value = "0.7211±0.0034"
$ws.Range("G5").Value = "0.5304±0.0044"
$ws.Range("B6").Value = "0.6438±0.0026"
$ws.Range("C6").Value = "0.5403±0.0046"
$ws.Range("D6").Value = "0.5034±0.0057"
$ws.Range("E6").Value = "0.2819±0.0068"
$ws.Range("F6").Value = "0.7017±0.0038"
$ws.Range("G6").Value = "0.5212±0.0049"
$ws.Range("B7").Value = "0.6622±0.0018"
$ws.Range("C7").Value = "0.6955±0.0104"
$ws.Range("D7").Value = "0.3979±0.0028"
$ws.Range("E7").Value = "0.2254±0.0074"
$ws.Range("F7").Value = "0.697±0.0024"
$ws.Range("G7").Value = "0.5062±0.0046"
$ws.Range("B8").Value = "0.6722±0.0018"
$ws.Range("C8").Value = "0.6336±0.0057"
$ws.Range("D8").Value = "0.4432±0.006"
$ws.Range("E8").Value = "0.2698±0.0066"
$ws.Range("F8").Value = "0.7041±0.0039"
$ws.Range("G8").Value = "0.5215±0.0039"
$ws.Range("B9").Value = "0.6536±0.011"
$ws.Range("C9").Value = "0.5599±0.0275"
$ws.Range("D9").Value = "0.5034±0.0153"
$ws.Range("E9").Value = "0.2924±0.0101"
$ws.Range("F9").Value = "0.7066±0.007"
$ws.Range("G9").Value = "0.5295±0.0072"
$ws.Range("B10").Value = "0.3463±0.1081"
$ws.Range("C10").Value = "0.4716±0.0288"
$ws.Range("D10").Value = "0.47±0.0125"
$ws.Range("E10").Value = "0.1812±0.039"
$ws.Range("F10").Value = "0.6258±0.0329"
$ws.Range("G10").Value = "0.4707±0.02"
$ws.Range("B11").Value = "0.576±0.0655"
$ws.Range("C11").Value = "0.4767±0.0284"
$ws.Range("D11").Value = "0.45±0.0197"
$ws.Range("E11").Value = "0.1925±0.0389"
$ws.Range("F11").Value = "0.5945±0.0359"
$ws.Range("G11").Value = "0.4629±0.0231"
$ws.Range("B12").Value = "0.6279±0.0122"
$ws.Range("C12").Value = "0.5338±0.0573"
$ws.Range("D12").Value = "0.4817±0.0308"
$ws.Range("E12").Value = "0.2514±0.0104"
$ws.Range("F12").Value = "0.6657±0.0107"
$ws.Range("G12").Value = "0.5033±0.0038"
$ws.Range("B13").Value = "0.686±0.0053"
$ws.Range("C13").Value = "0.6331±0.0053"
$ws.Range("D13").Value = "0.4925±0.0167"
$ws.Range("E13").Value = "0.3242±0.0193"
$ws.Range("F13").Value = "0.7385±0.0115"
$ws.Range("G13").Value = "0.5538±0.0116"
$ws.Range("B14").Value = "0.6835±0.0136"
$ws.Range("C14").Value = "0.6587±0.0417"
$ws.Range("D14").Value = "0.4728±0.0107"
$ws.Range("E14").Value = "0.3143±0.0133"
$ws.Range("F14").Value = "0.749±0.0157"
$ws.Range("G14").Value = "0.5496±0.011"

# Right table (TEST block): columns I-N, row by row
$ws.Range("I3").Value = "0.6856±0.003"
$ws.Range("J3").Value = "0.6454±0.0097"
$ws.Range("K3").Value = "0.4732±0.0032"
$ws.Range("L3").Value = "0.315±0.0079"
$ws.Range("M3").Value = "0.7131±0.0041"
$ws.Range("N3").Value = "0.546±0.0055"
$ws.Range("I5").Value = "0.6635±0.0035"
$ws.Range("J5").Value = "0.5688±0.006"
$ws.Range("K5").Value = "0.5141±0.0047"
$ws.Range("L5").Value = "0.3113±0.0072"
$ws.Range("M5").Value = "0.7325±0.0027"
$ws.Range("N5").Value = "0.54±0.0051"
$ws.Range("I6").Value = "0.6515±0.0043"
$ws.Range("J6").Value = "0.551±0.0075"
$ws.Range("K6").Value = "0.5069±0.0055"
$ws.Range("L6").Value = "0.2926±0.0085"
$ws.Range("M6").Value = "0.7093±0.0029"
$ws.Range("N6").Value = "0.528±0.0064"
$ws.Range("I7").Value = "0.6631±0.0014"
$ws.Range("J7").Value = "0.7037±0.0107"
$ws.Range("K7").Value = "0.3985±0.002"
$ws.Range("L7").Value = "0.2293±0.0059"
$ws.Range("M7").Value = "0.7037±0.0033"
$ws.Range("N7").Value = "0.5088±0.0041"
$ws.Range("I8").Value = "0.681±0.002"
$ws.Range("J8").Value = "0.6601±0.0085"
$ws.Range("K8").Value = "0.4555±0.0029"
$ws.Range("L8").Value = "0.2967±0.0061"
$ws.Range("M8").Value = "0.7258±0.0032"
$ws.Range("N8").Value = "0.539±0.0042"
$ws.Range("I9").Value = "0.6596±0.011"
$ws.Range("J9").Value = "0.5684±0.0292"
$ws.Range("K9").Value = "0.5059±0.0122"
$ws.Range("L9").Value = "0.3011±0.0116"
$ws.Range("M9").Value = "0.7117±0.005"
$ws.Range("N9").Value = "0.5347±0.0084"
$ws.Range("I10").Value = "0.3381±0.1162"
$ws.Range("J10").Value = "0.4827±0.0342"
$ws.Range("K10").Value = "0.4748±0.0162"
$ws.Range("L10").Value = "0.1887±0.0462"
$ws.Range("M10").Value = "0.6291±0.0387"
$ws.Range("N10").Value = "0.4786±0.0244"
$ws.Range("I11").Value = "0.5819±0.065"
$ws.Range("J11").Value = "0.4834±0.0332"
$ws.Range("K11").Value = "0.4553±0.0208"
$ws.Range("L11").Value = "0.2026±0.0429"
$ws.Range("M11").Value = "0.5992±0.0395"
$ws.Range("N11").Value = "0.4688±0.0261"
$ws.Range("I12").Value = "0.6353±0.0112"
$ws.Range("J12").Value = "0.5443±0.0612"
$ws.Range("K12").Value = "0.4854±0.0308"
$ws.Range("L12").Value = "0.2614±0.015"
$ws.Range("M12").Value = "0.671±0.013"
$ws.Range("N12").Value = "0.51±0.008"
$ws.Range("I13").Value = "0.6917±0.0042"
$ws.Range("J13").Value = "0.6454±0.0096"
$ws.Range("K13").Value = "0.5±0.0148"
$ws.Range("L13").Value = "0.3394±0.0153"
$ws.Range("M13").Value = "0.7475±0.0075"
$ws.Range("N13").Value = "0.5632±0.0081"
$ws.Range("I14").Value = "0.6879±0.0103"
$ws.Range("J14").Value = "0.674±0.0388"
$ws.Range("K14").Value = "0.475±0.014"
$ws.Range("L14").Value = "0.3246±0.0108"
$ws.Range("M14").Value = "0.7544±0.0151"
$ws.Range("N14").Value = "0.5563±0.0082"

# LN row (row 4) has no results for this run - clear its cells but keep
# them present as blank (default-styled) cells rather than removing them
$ws.Range("B4:G4").Value = ""
$ws.Range("B4:G4").Style = "Normal"
$ws.Range("I4:N4").Value = ""
$ws.Range("I4:N4").Style = "Normal"

# Update the active selection to reflect where the user left off
$ws.Range("L24").Select()
